# Update the division-problem worksheet table: replace the 25 problem
# strings (5 data rows x 5 columns) with the newly generated set of
# problems, preserving all formatting (font, size, paragraph settings).

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text  = "66÷5="
$tbl.Cell(1, 2).Range.Text  = "78÷7="
$tbl.Cell(1, 3).Range.Text  = "24÷8="
$tbl.Cell(1, 4).Range.Text  = "40÷4="
$tbl.Cell(1, 5).Range.Text  = "53÷7="

$tbl.Cell(5, 1).Range.Text  = "99÷6="
$tbl.Cell(5, 2).Range.Text  = "61÷6="
$tbl.Cell(5, 3).Range.Text  = "87÷8="
$tbl.Cell(5, 4).Range.Text  = "16÷4="
$tbl.Cell(5, 5).Range.Text  = "26÷5="

$tbl.Cell(9, 1).Range.Text  = "79÷9="
$tbl.Cell(9, 2).Range.Text  = "16÷3="
$tbl.Cell(9, 3).Range.Text  = "77÷7="
$tbl.Cell(9, 4).Range.Text  = "65÷2="
$tbl.Cell(9, 5).Range.Text  = "28÷8="

$tbl.Cell(13, 1).Range.Text = "60÷9="
$tbl.Cell(13, 2).Range.Text = "97÷6="
$tbl.Cell(13, 3).Range.Text = "37÷6="
$tbl.Cell(13, 4).Range.Text = "68÷9="
$tbl.Cell(13, 5).Range.Text = "30÷8="

$tbl.Cell(17, 1).Range.Text = "51÷5="
$tbl.Cell(17, 2).Range.Text = "77÷5="
$tbl.Cell(17, 3).Range.Text = "47÷8="
$tbl.Cell(17, 4).Range.Text = "18÷7="
$tbl.Cell(17, 5).Range.Text = "73÷5="
